$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (candidate 1)
$ws.Range("B2").Value = 80.27
$ws.Range("C2").Value = 'The candidate has a high score because they have experience working with the required technologies (ReactJS, JavaScript, CSS, Frontend Development, NextJS). Their projects showcase their proficiency in these skills, including building user-friendly interfaces, integrating language and computer vision models, and implementing research papers using PyTorch and Python. These experiences align well with the responsibilities of a Frontend Engineer Intern.'
$ws.Range("E2").Value = 'The applicant shows a strong willingness to come to Japan, citing the country''s leading position in AI and expressing interest in the work culture. Additionally, the applicant''s flexibility in working in different countries, as well as the willingness to adapt to new cultures, demonstrates a positive attitude towards international work opportunities. The career plan and expectations from the company are clear and aligned with the company''s focus. The applicant also displays strengths in teamwork, adaptability, and problem-solving, which are essential soft skills sought by the company. The applicant''s proactive approach to addressing weaknesses and the ability to derive valuable insights from complex data further demonstrate suitability for the role.'

# Row 3 (candidate 2)
$ws.Range("B3").Value = 74.54000000000001
$ws.Range("C3").Value = 'The applicant has a high score of 74.54 as their projects demonstrate strong skills in ReactJS, JavaScript, CSS, and frontend development, which are essential for the Frontend Engineer Intern position. The projects showcase their experience in building web applications, utilizing frameworks like Flutter and NodeJS, and integrating real-time functionalities. This indicates their ability to contribute efficiently to the team and create visually appealing web designs.'
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 'The applicant demonstrates a strong willingness to work in Japan due to his interest in AI, admiration for Japanese work culture, and willingness to learn a new language. The applicant also shows a good level of adaptability and openness to working in different countries. In terms of personality, the applicant displays strengths in teamwork, adapting to challenging situations, and having a unique background in engineering and graphic design. However, areas for improvement include the preference to not work alone and the need to acquire new skills more efficiently. Overall, the applicant''s positive attitude, willingness to learn, and alignment with the company''s soft skills make him a strong candidate for consideration.'

# Row 4 (candidate 3)
$ws.Range("B4").Value = 79.06999999999999
$ws.Range("C4").Value = 'The candidate has a moderate score of 79.07. This is because the candidate has experience in frontend development, specifically with ReactJS and CSS, which aligns with the job description. However, the candidate''s projects focus more on computer vision, machine learning, and web development using different frameworks, which may not directly relate to the required skills. This might have resulted in the moderate score.'
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 'The applicant shows a strong willingness to come to Japan, citing the country''s leadership in AI and the desire to be part of its work culture. They express openness to working in Japan after graduation and show a genuine interest in learning the language. Their career plan aligns with the company''s focus on AI/ML, and they have clear expectations for the company and the internship. In terms of personality, the applicant''s strengths in teamwork and adaptability, plus their ability to learn from setbacks and their unique background in both engineering and graphic design, demonstrate a well-rounded and adaptable individual. The only concern is their preference to not work alone, but their willingness to address this by trying to work alone shows potential for improvement.'

